# Apply the BOM update: add connector J4 designator to row 6,
# and bump its quantity in row 8 (0.1" Pin Header, 01x08) from 1 to 2.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 (0.1" Pin Header / 01x08): schematic designator becomes "J4"
$ws.Range("C6").Value = "J4"

# Row 8 (SMD Resistor 47 kOhm, R2/R3): quantity per board 1 -> 2
# (H8 = G8*F8 recalculates automatically from 0.1 to 0.2)
$ws.Range("F8").Value = 2

# Update the active selection left by the editor
$ws.Range("C7").Select()
